# Generate Report for Handback
# The handback transform for file "8f4f9e54-ae52-456d-ae8b-1312b1766814" failed
# because the handback file name did not match the handoff file name. Update the
# status for that row on all three sheets and record the error detail for the
# zh-cn and de-de language rows.

$wb = $excel.ActiveWorkbook
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: both the zh-cn and de-de status columns for this file
$overview.Range("B7").Value = $newStatus
$overview.Range("C7").Value = $newStatus

# Per-language sheets: Status column (C) for the same file row
$zhcn.Range("C7").Value = $newStatus
$dede.Range("C7").Value = $newStatus

# Per-language sheets: Error Detail column (L) explaining the failure
$zhcn.Range("L7").Value = "Handback file name: ekqtdr1b.rna is different with handoff file name: 8f4f9e54-ae52-456d-ae8b-1312b1766814.584dbb32665d7b0a821170b106ee983b6c80ee93.zh-cn."
$dede.Range("L7").Value = "Handback file name: ekqtdr1b.rna is different with handoff file name: 8f4f9e54-ae52-456d-ae8b-1312b1766814.584dbb32665d7b0a821170b106ee983b6c80ee93.de-de."
